$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Fill in "nan" placeholders for the existing row 22 blank cells (B:K and M)
$ws.Range("B22").Value = "nan"
$ws.Range("C22").Value = "nan"
$ws.Range("D22").Value = "nan"
$ws.Range("E22").Value = "nan"
$ws.Range("F22").Value = "nan"
$ws.Range("G22").Value = "nan"
$ws.Range("H22").Value = "nan"
$ws.Range("I22").Value = "nan"
$ws.Range("J22").Value = "nan"
$ws.Range("K22").Value = "nan"
$ws.Range("M22").Value = "nan"

# Add the new event row (row 23)
$ws.Range("A23").Value = "'20"
$ws.Range("A23").Style = "Normal"
$ws.Range("L23").Value = "12\5\2025"
$ws.Range("N23").Value = "تم سن السليندر وتغيير الجرئد الاماميه (1_2_4_5_7_8)"
$ws.Range("O23").Value = "الخبير"
